$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-29 Friday", "2024-11-30 Saturday"),
    @("60×21=", "55×26="),
    @("47×35=", "56×77="),
    @("91×57=", "38×98="),
    @("25×46=", "29×26="),
    @("37×73=", "39×56="),
    @("12×22=", "29×17="),
    @("97×92=", "19×42="),
    @("41×52=", "15×50="),
    @("85×75=", "57×90="),
    @("85×84=", "68×66="),
    @("81×22=", "58×75="),
    @("20×70=", "11×77="),
    @("29×55=", "84×99="),
    @("21×43=", "75×53="),
    @("35×43=", "79×13="),
    @("87×35=", "88×89="),
    @("39×16=", "65×98="),
    @("95×80=", "94×15="),
    @("80×64=", "38×81="),
    @("70×61=", "39×87="),
    @("43×14=", "34×41="),
    @("52×54=", "80×97="),
    @("21×84=", "86×20="),
    @("38×61=", "18×33="),
    @("40×47=", "18×84=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
